# Weekly driver report update for 2025-04-19
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Bad Drivers" table -------------------------------------------------
# Update the existing bad-driver row (A3:D3) with refreshed figures.
$ws.Range("A3").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.40.0.4"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 202
$ws.Range("D3").Value = 96.90000000000001

# Insert a brand-new bad-driver row beneath it (pushes Totals: + everything
# below down by one row, inheriting the row-3 number formatting).
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.40.0.7"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 62
$ws.Range("D4").Value = 98.2

# Totals: row (now row 5 after the insert) - refresh the summed figures.
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 264

# --- "Good Drivers" table --------------------------------------------------
# The table itself (header now at row 12, data rows 13-18) kept its 6 rows
# of data, just reordered with refreshed sample counts; two rows lost their
# "Driver Vintage" date. Overwrite each data row in place.

$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B13").Value = 56018
$ws.Range("D13").Value = 100
$ws.Range("E13").ClearContents()

$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B14").Value = 34244
$ws.Range("D14").Value = 100
$ws.Range("E14").ClearContents()

$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B15").Value = 442178
$ws.Range("D15").Value = 99.90000000000001
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2024-11-10"
$ws.Range("D15").Copy()
$ws.Range("E15").PasteSpecial(-4122)

$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B16").Value = 77849
$ws.Range("D16").Value = 99.90000000000001
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2021-08-18"
$ws.Range("D16").Copy()
$ws.Range("E16").PasteSpecial(-4122)

$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B17").Value = 59673
$ws.Range("D17").Value = 100
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2020-08-05"
$ws.Range("D17").Copy()
$ws.Range("E17").PasteSpecial(-4122)

$ws.Range("A18").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B18").Value = 113652
$ws.Range("D18").Value = 100
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2019-12-14"
$ws.Range("D18").Copy()
$ws.Range("E18").PasteSpecial(-4122)
